$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Diákok"
$ws.Range("A3").Value = "Termek"
$ws.Range("A4").Value = "Tanári szobák"
$ws.Range("A5").Value = "Interaktív táblák"
$ws.Range("A6").Value = "Nyomtatók"

$ws.Range("A6").Select()
